# Generate Report for Handoff
# - Updates status text from "Handed back: in sync with en-US" to "Ready for handoff"
# - Refreshes the handoff/generate timestamps
# - Narrows the "Status"/"zh-cn"/"de-de" column widths

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2017-02-21 10:50:14"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2017-02-21 10:49:56"

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2017-02-21 10:50:14"

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
